$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# area4 (column E) and area5 (column F) value corrections for rows 2-4
$ws.Range("E2").Value = 208
$ws.Range("F2").Value = 84

$ws.Range("E3").Value = 172

$ws.Range("E4").Value = 86
$ws.Range("F4").Value = 39
